$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateLabel = "Sunday, Jan 15"

# Row 169 (NUMBER=168)
$ws.Cells.Item(169, 1).Value = 168
$ws.Cells.Item(169, 2).Value = $dateLabel
$ws.Cells.Item(169, 3).Value = "6:10 AM"
$ws.Cells.Item(169, 4).Value = "FR2022"
$ws.Cells.Item(169, 5).Value = "Dublin"
$ws.Cells.Item(169, 6).Value = "(DUB)"
$ws.Cells.Item(169, 7).Value = "Buzz "
$ws.Cells.Item(169, 8).Value = "B38M"
$ws.Cells.Item(169, 9).Value = "(SP-RZG)"
$ws.Cells.Item(169, 10).Value = "6:16 AM"
$ws.Cells.Item(169, 11).Borders.LineStyle = -4142
$ws.Cells.Item(169, 12).Value = "0 hours, 6 minutes"
$ws.Cells.Item(169, 13).Borders.LineStyle = -4142

# Row 170 (NUMBER=169)
$ws.Cells.Item(170, 1).Value = 169
$ws.Cells.Item(170, 2).Value = $dateLabel
$ws.Cells.Item(170, 3).Value = "6:15 AM"
$ws.Cells.Item(170, 4).Value = "FR1894"
$ws.Cells.Item(170, 5).Value = "Amman"
$ws.Cells.Item(170, 6).Value = "(AMM)"
$ws.Cells.Item(170, 7).Value = "Ryanair "
$ws.Cells.Item(170, 8).Value = "B738"
$ws.Cells.Item(170, 9).Value = "(SP-RSV)"
$ws.Cells.Item(170, 10).Value = "6:20 AM"
$ws.Cells.Item(170, 11).Borders.LineStyle = -4142
$ws.Cells.Item(170, 12).Value = "0 hours, 5 minutes"
$ws.Cells.Item(170, 13).Borders.LineStyle = -4142

# Row 171 (NUMBER=170)
$ws.Cells.Item(171, 1).Value = 170
$ws.Cells.Item(171, 2).Value = $dateLabel
$ws.Cells.Item(171, 3).Value = "7:10 AM"
$ws.Cells.Item(171, 4).Value = "FR2350"
$ws.Cells.Item(171, 5).Value = "Shannon"
$ws.Cells.Item(171, 6).Value = "(SNN)"
$ws.Cells.Item(171, 7).Value = "Ryanair "
$ws.Cells.Item(171, 8).Value = "B38M"
$ws.Cells.Item(171, 9).Value = "(SP-RZO)"
$ws.Cells.Item(171, 10).Value = "7:13 AM"
$ws.Cells.Item(171, 11).Borders.LineStyle = -4142
$ws.Cells.Item(171, 12).Value = "0 hours, 3 minutes"
$ws.Cells.Item(171, 13).Borders.LineStyle = -4142

# Row 172 (NUMBER=171)
$ws.Cells.Item(172, 1).Value = 171
$ws.Cells.Item(172, 2).Value = $dateLabel
$ws.Cells.Item(172, 3).Value = "7:40 AM"
$ws.Cells.Item(172, 4).Value = "FR4528"
$ws.Cells.Item(172, 5).Value = "Oslo"
$ws.Cells.Item(172, 6).Value = "(TRF)"
$ws.Cells.Item(172, 7).Value = "Ryanair "
$ws.Cells.Item(172, 8).Value = "B738"
$ws.Cells.Item(172, 9).Value = "(SP-RKP)"
$ws.Cells.Item(172, 10).Value = "7:44 AM"
$ws.Cells.Item(172, 11).Borders.LineStyle = -4142
$ws.Cells.Item(172, 12).Value = "0 hours, 4 minutes"
$ws.Cells.Item(172, 13).Borders.LineStyle = -4142

# Row 173 (NUMBER=172)
$ws.Cells.Item(173, 1).Value = 172
$ws.Cells.Item(173, 2).Value = $dateLabel
$ws.Cells.Item(173, 3).Value = "7:40 AM"
$ws.Cells.Item(173, 4).Value = "FR9259"
$ws.Cells.Item(173, 5).Value = "Malta"
$ws.Cells.Item(173, 6).Value = "(MLA)"
$ws.Cells.Item(173, 7).Value = "Ryanair "
$ws.Cells.Item(173, 8).Value = "B738"
$ws.Cells.Item(173, 9).Value = "(SP-RSP)"
$ws.Cells.Item(173, 10).Value = "7:39 AM"
$ws.Cells.Item(173, 11).Borders.LineStyle = -4142
$ws.Cells.Item(173, 12).Value = "0 hours, -1 minutes"
$ws.Cells.Item(173, 13).Borders.LineStyle = -4142

# Row 174 (NUMBER=173)
$ws.Cells.Item(174, 1).Value = 173
$ws.Cells.Item(174, 2).Value = $dateLabel
$ws.Cells.Item(174, 3).Value = "8:25 AM"
$ws.Cells.Item(174, 4).Value = "FR4060"
$ws.Cells.Item(174, 5).Value = "Malaga"
$ws.Cells.Item(174, 6).Value = "(AGP)"
$ws.Cells.Item(174, 7).Value = "Buzz "
$ws.Cells.Item(174, 8).Value = "B38M"
$ws.Cells.Item(174, 9).Value = "(SP-RZE)"
$ws.Cells.Item(174, 10).Value = "8:34 AM"
$ws.Cells.Item(174, 11).Borders.LineStyle = -4142
$ws.Cells.Item(174, 12).Value = "0 hours, 9 minutes"
$ws.Cells.Item(174, 13).Borders.LineStyle = -4142

# Row 175 (NUMBER=174)
$ws.Cells.Item(175, 1).Value = 174
$ws.Cells.Item(175, 2).Value = $dateLabel
$ws.Cells.Item(175, 3).Value = "9:50 AM"
$ws.Cells.Item(175, 4).Value = "FR1950"
$ws.Cells.Item(175, 5).Value = "Manchester"
$ws.Cells.Item(175, 6).Value = "(MAN)"
$ws.Cells.Item(175, 7).Value = "Ryanair "
$ws.Cells.Item(175, 8).Value = "B738"
$ws.Cells.Item(175, 9).Value = "(EI-EBW)"
$ws.Cells.Item(175, 10).Value = "9:50 AM"
$ws.Cells.Item(175, 11).Borders.LineStyle = -4142
$ws.Cells.Item(175, 12).Value = "0 hours, 0 minutes"
$ws.Cells.Item(175, 13).Borders.LineStyle = -4142

# Row 176 (NUMBER=175)
$ws.Cells.Item(176, 1).Value = 175
$ws.Cells.Item(176, 2).Value = $dateLabel
$ws.Cells.Item(176, 3).Value = "10:00 AM"
$ws.Cells.Item(176, 4).Value = "FR2008"
$ws.Cells.Item(176, 5).Value = "London"
$ws.Cells.Item(176, 6).Value = "(STN)"
$ws.Cells.Item(176, 7).Value = "Ryanair "
$ws.Cells.Item(176, 8).Value = "B38M"
$ws.Cells.Item(176, 9).Value = "(EI-HAX)"
$ws.Cells.Item(176, 10).Value = "10:11 AM"
$ws.Cells.Item(176, 11).Borders.LineStyle = -4142
$ws.Cells.Item(176, 12).Value = "0 hours, 11 minutes"
$ws.Cells.Item(176, 13).Borders.LineStyle = -4142

# Row 177 (NUMBER=176)
$ws.Cells.Item(177, 1).Value = 176
$ws.Cells.Item(177, 2).Value = $dateLabel
$ws.Cells.Item(177, 3).Value = "11:30 AM"
$ws.Cells.Item(177, 4).Value = "FR8084"
$ws.Cells.Item(177, 5).Value = "Birmingham"
$ws.Cells.Item(177, 6).Value = "(BHX)"
$ws.Cells.Item(177, 7).Value = "Ryanair "
$ws.Cells.Item(177, 8).Value = "B738"
$ws.Cells.Item(177, 9).Value = "(EI-EXD)"
$ws.Cells.Item(177, 10).Value = "11:41 AM"
$ws.Cells.Item(177, 11).Borders.LineStyle = -4142
$ws.Cells.Item(177, 12).Value = "0 hours, 11 minutes"
$ws.Cells.Item(177, 13).Borders.LineStyle = -4142

# Row 178 (NUMBER=177)
$ws.Cells.Item(178, 1).Value = 177
$ws.Cells.Item(178, 2).Value = $dateLabel
$ws.Cells.Item(178, 3).Value = "1:25 PM"
$ws.Cells.Item(178, 4).Value = "FR1888"
$ws.Cells.Item(178, 5).Value = "Paris"
$ws.Cells.Item(178, 6).Value = "(BVA)"
$ws.Cells.Item(178, 7).Value = "Buzz "
$ws.Cells.Item(178, 8).Value = "B38M"
$ws.Cells.Item(178, 9).Value = "(SP-RZG)"
$ws.Cells.Item(178, 10).Value = "1:32 PM"
$ws.Cells.Item(178, 11).Borders.LineStyle = -4142
$ws.Cells.Item(178, 12).Value = "0 hours, 7 minutes"
$ws.Cells.Item(178, 13).Borders.LineStyle = -4142

# Row 179 (NUMBER=178)
$ws.Cells.Item(179, 1).Value = 178
$ws.Cells.Item(179, 2).Value = $dateLabel
$ws.Cells.Item(179, 3).Value = "2:40 PM"
$ws.Cells.Item(179, 4).Value = "FR1942"
$ws.Cells.Item(179, 5).Value = "Bologna"
$ws.Cells.Item(179, 6).Value = "(BLQ)"
$ws.Cells.Item(179, 7).Value = "Ryanair "
$ws.Cells.Item(179, 8).Value = "B738"
$ws.Cells.Item(179, 9).Value = "(9H-QDA)"
$ws.Cells.Item(179, 10).Value = "2:40 PM"
$ws.Cells.Item(179, 11).Borders.LineStyle = -4142
$ws.Cells.Item(179, 12).Value = "0 hours, 0 minutes"
$ws.Cells.Item(179, 13).Borders.LineStyle = -4142

# Row 180 (NUMBER=179)
$ws.Cells.Item(180, 1).Value = 179
$ws.Cells.Item(180, 2).Value = $dateLabel
$ws.Cells.Item(180, 3).Value = "3:15 PM"
$ws.Cells.Item(180, 4).Value = "FR1934"
$ws.Cells.Item(180, 5).Value = "Helsinki"
$ws.Cells.Item(180, 6).Value = "(HEL)"
$ws.Cells.Item(180, 7).Value = "Ryanair "
$ws.Cells.Item(180, 8).Value = "B738"
$ws.Cells.Item(180, 9).Value = "(SP-RKD)"
$ws.Cells.Item(180, 10).Value = "3:11 PM"
$ws.Cells.Item(180, 11).Borders.LineStyle = -4142
$ws.Cells.Item(180, 12).Value = "0 hours, -4 minutes"
$ws.Cells.Item(180, 13).Borders.LineStyle = -4142

# Row 181 (NUMBER=180)
$ws.Cells.Item(181, 1).Value = 180
$ws.Cells.Item(181, 2).Value = $dateLabel
$ws.Cells.Item(181, 3).Value = "3:35 PM"
$ws.Cells.Item(181, 4).Value = "FR1898"
$ws.Cells.Item(181, 5).Value = "Eindhoven"
$ws.Cells.Item(181, 6).Value = "(EIN)"
$ws.Cells.Item(181, 7).Value = "Ryanair "
$ws.Cells.Item(181, 8).Value = "B738"
$ws.Cells.Item(181, 9).Value = "(SP-RSV)"
$ws.Cells.Item(181, 10).Value = "3:35 PM"
$ws.Cells.Item(181, 11).Borders.LineStyle = -4142
$ws.Cells.Item(181, 12).Value = "0 hours, 0 minutes"
$ws.Cells.Item(181, 13).Borders.LineStyle = -4142
